$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Bill of Materials table ---
# Two rows were removed from the original table ("1uF" capacitor and "1k"
# resistor), so every row below them shifts up by one/two places
# respectively.  Rather than physically deleting rows (which would also
# shrink the sheet's used range / dimension and disturb the untouched
# trailing blank rows 31-32), we simply rewrite the cell values for the
# affected range A3:C28 to their new post-edit contents, and blank out
# what used to be the last two data rows (29-30).

$data = @(
    @("Capacitor 0805", "0.001uF", $null),
    @("Capacitor 0805", "0.1uF", $null),
    @("Capacitor 0805", "2.2uF", $null),
    @("Capacitor 0805", "4.7uF", $null),
    @("Capacitor 0805", "47pF", $null),
    @("Resistor 0805", "33R", $null),
    @("Resistor 0805", "330R", $null),
    @("Resistor 0805", "1.5k", $null),
    @("Resistor 0805", "2k", $null),
    @("Resistor 0805", "10k", $null),
    @("Inductor", "2.2uH", "810-MLZ2012A2R2W"),
    @("LED", "Red", "720-LHR974-LP-1"),
    @("Transistor", "PMOS", "771-PMV160UP215"),
    @("Transistor", "NMOS", "771-NX3008NBK,215"),
    @("Switch", "Toggle", "688-SSSS810701"),
    @("Switch", "Push", "612-TL1015AF"),
    @("Connector", "0.05`" 2x5", "649-221111-00010T4LF"),
    @("Connector", "uUSB", "649-10118192-0001LF"),
    @("Connector", "LiPo", "455-1749-1-ND (DigiKey)"),
    @("Header", "Female 1x10", "992-10FX1-254MM"),
    @("Header", "8x 2 pin jump", "649-68000-236HLF"),
    @("Regulator", "LP5907 3.3V", "LP5907MFX-3.3/NOPB"),
    @("Regulator", "TPS62730", "595-TPS62730DRYT"),
    @("Charger", "MCP73831", "579-MCP73831T-2ACIOT"),
    @("Bluetooth", "BLE112", "603-BLE112-A"),
    @("Jumper", "2 Pin", "806-SX1100-B")
)

$row = 3
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    if ($entry[2] -eq $null) {
        $ws.Cells.Item($row, 3).Value = $null
    } else {
        $ws.Cells.Item($row, 3).Value = $entry[2]
    }
    $row = $row + 1
}

# Former rows 29-30 no longer hold BOM entries now that the table shrank.
$ws.Range("A29:C30").Value = $null

# The "NMOS" Mouser-number cell (now row 18) keeps the word-wrap styling
# that used to belong to the "Push" switch row, while that row (now row
# 20) reverts to the plain/default style.
$ws.Range("C20").Style = "Normal"
$ws.Range("C18").WrapText = $true

# --- Misc bookkeeping to match the saved workbook state ---
$ws.Range("H10").Select()
